# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.683.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.699.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "677.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.497"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.703.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.735.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "472.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.851.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.691.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0908"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "169.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.268"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.64%  "
